# Auto-generated edit script applying the committed numeric updates
# to the Sagittarius_Profits workbook (profit-calc columns H-N) across
# all eight sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4423.3335
$ws.Range("J17").Value = 4423.3335
$ws.Range("L17").Value = 13270.0005
$ws.Range("N17").Value = -13606.0005
$ws.Range("H19").Value = 1335.25
$ws.Range("I19").Value = 1030.3334
$ws.Range("J19").Value = 2250
$ws.Range("K19").Value = 1030.3334
$ws.Range("L19").Value = 2250
$ws.Range("M19").Value = -855.3334
$ws.Range("N19").Value = -2600
$ws.Range("H28").Value = 1079.6923
$ws.Range("I28").Value = 1079.6923
$ws.Range("K28").Value = 1079.6923
$ws.Range("M28").Value = -594.6922999999999
$ws.Range("H62").Value = 7934.4287
$ws.Range("I62").Value = 7889
$ws.Range("K62").Value = 7889
$ws.Range("M62").Value = -7265
$ws.Range("H65").Value = 7934.4287
$ws.Range("I65").Value = 7889
$ws.Range("K65").Value = 39445
$ws.Range("M65").Value = -36325
$ws.Range("H98").Value = 3438
$ws.Range("J98").Value = 4000
$ws.Range("L98").Value = 4000
$ws.Range("N98").Value = -6996
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H122").Value = 3438
$ws.Range("J122").Value = 4000
$ws.Range("L122").Value = 12000
$ws.Range("N122").Value = -16900
$ws.Range("H132").Value = 1049.8125
$ws.Range("I132").Value = 985.9286
$ws.Range("K132").Value = 2957.7858
$ws.Range("M132").Value = -427.7857999999997
$ws.Range("H135").Value = 771.375
$ws.Range("I135").Value = 771.375
$ws.Range("K135").Value = 6942.375
$ws.Range("M135").Value = -4407.375
$ws.Range("H137").Value = 2062.1667
$ws.Range("I137").Value = 2062.1667
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 6186.500100000001
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -3636.500100000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3404.1155
$ws.Range("I32").Value = 3404.1155
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3404.1155
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -3117.1155
$ws.Range("H46").Value = 1941.5
$ws.Range("J46").Value = 1941.5
$ws.Range("L46").Value = 1941.5
$ws.Range("N46").Value = -2579.5
$ws.Range("H61").Value = 4049.875
$ws.Range("I61").Value = 1884.5454
$ws.Range("J61").Value = 5882.077
$ws.Range("K61").Value = 1884.5454
$ws.Range("L61").Value = 5882.077
$ws.Range("M61").Value = -1672.5454
$ws.Range("N61").Value = -6306.077
$ws.Range("H74").Value = 1781.5834
$ws.Range("I74").Value = 1675.6111
$ws.Range("J74").Value = 2099.5
$ws.Range("K74").Value = 1675.6111
$ws.Range("L74").Value = 2099.5
$ws.Range("M74").Value = -801.6111000000001
$ws.Range("N74").Value = -3847.5
$ws.Range("H77").Value = 1781.5834
$ws.Range("I77").Value = 1675.6111
$ws.Range("J77").Value = 2099.5
$ws.Range("K77").Value = 8378.0555
$ws.Range("L77").Value = 10497.5
$ws.Range("M77").Value = -4010.0555
$ws.Range("N77").Value = -19233.5
$ws.Range("H102").Value = 1250
$ws.Range("I102").Value = 1250
$ws.Range("K102").Value = 1250
$ws.Range("M102").Value = 372
$ws.Range("H113").Value = 130499.664
$ws.Range("J113").Value = 130499.664
$ws.Range("L113").Value = 130499.664
$ws.Range("N113").Value = -139177.664
$ws.Range("H132").Value = 3012
$ws.Range("I132").Value = 3012
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9036
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -6506
$ws.Range("H136").Value = 4049.875
$ws.Range("I136").Value = 1884.5454
$ws.Range("J136").Value = 5882.077
$ws.Range("K136").Value = 5653.6362
$ws.Range("L136").Value = 17646.231
$ws.Range("M136").Value = -3103.6362
$ws.Range("N136").Value = -22746.231

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 7999
$ws.Range("I82").Value = 7999
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 7999
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -7616
$ws.Range("H85").Value = 7999
$ws.Range("I85").Value = 7999
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 7999
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -6673
$ws.Range("H99").Value = 2599.5
$ws.Range("I99").Value = 2219.4
$ws.Range("K99").Value = 2219.4
$ws.Range("M99").Value = -721.4000000000001
$ws.Range("H105").Value = 1950
$ws.Range("I105").Value = 1950
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1950
$ws.Range("L105").ClearContents()
$ws.Range("N105").Value = 0
$ws.Range("M105").Value = -203
$ws.Range("H134").Value = 1595.4
$ws.Range("I134").Value = 1595.4
$ws.Range("K134").Value = 4786.200000000001
$ws.Range("M134").Value = -2251.200000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 9523.5
$ws.Range("I103").Value = 9523.5
$ws.Range("K103").Value = 9523.5
$ws.Range("M103").Value = -8351.5
$ws.Range("I122").Value = 12
$ws.Range("J122").Value = 1703.5
$ws.Range("K122").Value = 36
$ws.Range("L122").Value = 5110.5
$ws.Range("M122").Value = 2414
$ws.Range("N122").Value = -10010.5
$ws.Range("H130").Value = 48593
$ws.Range("J130").Value = 48593
$ws.Range("L130").Value = 48593
$ws.Range("N130").Value = -58633
$ws.Range("H132").Value = 1463.08
$ws.Range("I132").Value = 1441.238
$ws.Range("J132").Value = 1577.75
$ws.Range("K132").Value = 4323.714
$ws.Range("L132").Value = 4733.25
$ws.Range("M132").Value = -1793.714
$ws.Range("N132").Value = -9793.25
$ws.Range("H134").Value = 1868.8422
$ws.Range("I134").Value = 1816.2354
$ws.Range("K134").Value = 5448.706200000001
$ws.Range("M134").Value = -2913.706200000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 208.28572
$ws.Range("I7").Value = 218
$ws.Range("K7").Value = 654
$ws.Range("M7").Value = -542
$ws.Range("H34").Value = 340
$ws.Range("I34").Value = 340
$ws.Range("K34").Value = 1020
$ws.Range("M34").Value = -936
$ws.Range("H38").Value = 53
$ws.Range("I38").Value = 43.6
$ws.Range("J38").Value = 100
$ws.Range("K38").Value = 130.8
$ws.Range("L38").Value = 300
$ws.Range("M38").Value = 216.2
$ws.Range("N38").Value = -994
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").ClearContents()
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = 0
$ws.Range("H55").Value = 1738
$ws.Range("H92").Value = 240.6
$ws.Range("I92").Value = 166.26086
$ws.Range("K92").Value = 498.7825800000001
$ws.Range("M92").Value = 749.2174199999999
$ws.Range("H139").Value = 1410.3334
$ws.Range("I139").Value = 1410.3334
$ws.Range("K139").Value = 4231.0002
$ws.Range("M139").Value = 908.9997999999996
$ws.Range("H140").Value = 5687.1577
$ws.Range("I140").Value = 1389.1538
$ws.Range("K140").Value = 4167.4614
$ws.Range("M140").Value = 1012.5386

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8493.9375
$ws.Range("I70").Value = 8543.571
$ws.Range("K70").Value = 8543.571
$ws.Range("M70").Value = -8273.571
$ws.Range("H73").Value = 8493.9375
$ws.Range("I73").Value = 8543.571
$ws.Range("K73").Value = 8543.571
$ws.Range("M73").Value = -7607.571
$ws.Range("H102").Value = 1327.2222
$ws.Range("I102").Value = 1327.2222
$ws.Range("K102").Value = 1327.2222
$ws.Range("M102").Value = 294.7778000000001
$ws.Range("H113").Value = 1720
$ws.Range("I113").Value = 1720
$ws.Range("K113").Value = 1720
$ws.Range("M113").Value = 450
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").ClearContents()
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = 0

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3790.7693
$ws.Range("I40").Value = 3111.7144
$ws.Range("K40").Value = 3111.7144
$ws.Range("M40").Value = -2975.7144
$ws.Range("H61").Value = 5972.25
$ws.Range("I61").Value = 5972.25
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 5972.25
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -5770.25
$ws.Range("H100").Value = 4714.2856
$ws.Range("J100").Value = 4714.2856
$ws.Range("L100").Value = 4714.2856
$ws.Range("N100").Value = -5796.2856
$ws.Range("H113").Value = 5972.25
$ws.Range("I113").Value = 5972.25
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 5972.25
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -3802.25
$ws.Range("H132").Value = 3527.5557
$ws.Range("I132").Value = 3257.1538
$ws.Range("J132").Value = 4230.6
$ws.Range("K132").Value = 9771.4614
$ws.Range("L132").Value = 12691.8
$ws.Range("M132").Value = -7241.4614
$ws.Range("N132").Value = -17751.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 11113880
$ws.Range("I100").Value = 12501865
$ws.Range("K100").Value = 25003730
$ws.Range("M100").Value = -25003189
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 0
$ws.Range("H132").Value = 718.6
$ws.Range("I132").Value = 698.25
$ws.Range("K132").Value = 2094.75
$ws.Range("M132").Value = 435.25
